$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 570.7143
$ws.Range("I6").Value = 615.8333
$ws.Range("K6").Value = 1847.4999
$ws.Range("M6").Value = -1735.4999

# Row 12
$ws.Range("H12").Value = 799.2
$ws.Range("I12").Value = 899
$ws.Range("K12").Value = 899
$ws.Range("M12").Value = -729

# Row 17
$ws.Range("H17").Value = 1536.3243
$ws.Range("J17").Value = 1536.3243
$ws.Range("L17").Value = 4608.9729
$ws.Range("N17").Value = -4944.9729

# Row 33
$ws.Range("H33").Value = 311.4
$ws.Range("I33").Value = 344.17648
$ws.Range("K33").Value = 344.17648
$ws.Range("M33").Value = -115.17648

# Row 40
$ws.Range("H40").Value = 1915.375
$ws.Range("I40").Value = 1293.75
$ws.Range("J40").Value = 2537
$ws.Range("K40").Value = 1293.75
$ws.Range("L40").Value = 2537
$ws.Range("M40").Value = -1118.75
$ws.Range("N40").Value = -2887

# Row 70
$ws.Range("H70").Value = 4675.25
$ws.Range("J70").Value = 1800
$ws.Range("L70").Value = 5400
$ws.Range("N70").Value = -5940

# Row 73
$ws.Range("H73").Value = 4675.25
$ws.Range("J73").Value = 1800
$ws.Range("L73").Value = 5400
$ws.Range("N73").Value = -7272

# Row 106
$ws.Range("H106").Value = 14119.066
$ws.Range("I106").Value = 11368.385
$ws.Range("K106").Value = 11368.385
$ws.Range("M106").Value = -10737.385

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 978.7368
$ws.Range("I2").Value = 1018.05884
$ws.Range("K2").Value = 1018.05884
$ws.Range("M2").Value = -905.05884

# Row 32
$ws.Range("H32").Value = 12384.45
$ws.Range("I32").Value = 12384.45
$ws.Range("K32").Value = 12384.45
$ws.Range("M32").Value = -12097.45

# Row 45
$ws.Range("H45").Value = 2110.5715
$ws.Range("I45").Value = 2145.6667
$ws.Range("K45").Value = 2145.6667
$ws.Range("M45").Value = -1768.6667

# Row 61
$ws.Range("H61").Value = 1462.3334
$ws.Range("I61").Value = 1462.3334
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1462.3334
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -1250.3334

# Row 63
$ws.Range("H63").Value = 5800.8335
$ws.Range("I63").Value = 2279.6785
$ws.Range("K63").Value = 2279.6785
$ws.Range("M63").Value = -1593.6785

# Row 66
$ws.Range("H66").Value = 5800.8335
$ws.Range("I66").Value = 2279.6785
$ws.Range("K66").Value = 11398.3925
$ws.Range("M66").Value = -7966.3925

# Row 110
$ws.Range("H110").Value = 4028.1875
$ws.Range("I110").Value = 4256.1333
$ws.Range("K110").Value = 4256.1333
$ws.Range("M110").Value = -2211.1333

# Row 116
$ws.Range("H116").Value = 978.7368
$ws.Range("I116").Value = 1018.05884
$ws.Range("K116").Value = 1018.05884
$ws.Range("M116").Value = 1275.94116

# Row 122
$ws.Range("H122").Value = 2424.8
$ws.Range("I122").Value = 2424.8
$ws.Range("K122").Value = 7274.400000000001
$ws.Range("M122").Value = -4824.400000000001

# Row 132
$ws.Range("H132").Value = 3230.7896
$ws.Range("I132").Value = 1499.25
$ws.Range("J132").Value = 6199.143
$ws.Range("K132").Value = 4497.75
$ws.Range("L132").Value = 18597.429
$ws.Range("M132").Value = -1967.75
$ws.Range("N132").Value = -23657.429

# Row 136
$ws.Range("H136").Value = 1462.3334
$ws.Range("I136").Value = 1462.3334
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4387.0002
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -1837.0002

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 978.7368
$ws.Range("I3").Value = 1018.05884
$ws.Range("K3").Value = 1018.05884
$ws.Range("M3").Value = -904.05884

# Row 99
$ws.Range("H99").Value = 1493.75
$ws.Range("I99").Value = 1493.75
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1493.75
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 4.25

# Row 134
$ws.Range("H134").Value = 3083.25
$ws.Range("I134").Value = 3083.25
$ws.Range("K134").Value = 9249.75
$ws.Range("M134").Value = -6714.75

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1833
$ws.Range("I16").Value = 1833
$ws.Range("K16").Value = 1833
$ws.Range("M16").Value = -1546

# Row 94
$ws.Range("H94").Value = 322.8

# Row 99
$ws.Range("H99").Value = 2298.75
$ws.Range("J99").Value = 2450
$ws.Range("L99").Value = 2450
$ws.Range("N99").Value = -5446

# Row 107
$ws.Range("H107").Value = 674
$ws.Range("I107").Value = 620
$ws.Range("J107").Value = 998
$ws.Range("K107").Value = 620
$ws.Range("L107").Value = 998
$ws.Range("M107").Value = 1300
$ws.Range("N107").Value = -4838

# Row 113
$ws.Range("H113").Value = 1833
$ws.Range("I113").Value = 1833
$ws.Range("K113").Value = 1833
$ws.Range("M113").Value = 337

# Row 126
$ws.Range("H126").Value = 2298.75
$ws.Range("J126").Value = 2450
$ws.Range("L126").Value = 7350
$ws.Range("N126").Value = -12290

# Row 132
$ws.Range("H132").Value = 2406.4
$ws.Range("I132").Value = 2406.4
$ws.Range("K132").Value = 7219.200000000001
$ws.Range("M132").Value = -4689.200000000001

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1164.3334
$ws.Range("J5").Value = 1247
$ws.Range("L5").Value = 3741
$ws.Range("N5").Value = -3965

# Row 15
$ws.Range("H15").Value = 425.42856
$ws.Range("I15").Value = 230.5
$ws.Range("K15").Value = 691.5
$ws.Range("M15").Value = -551.5

# Row 23
$ws.Range("H23").Value = 508.5
$ws.Range("I23").Value = 512.6667
$ws.Range("K23").Value = 1538.0001
$ws.Range("M23").Value = -1303.0001

# Row 135
$ws.Range("H135").Value = 1164.3334
$ws.Range("J135").Value = 1247
$ws.Range("L135").Value = 11223
$ws.Range("N135").Value = -16293

# Row 139
$ws.Range("H139").Value = 2891.25
$ws.Range("I139").Value = 2756.7144
$ws.Range("J139").Value = 3833
$ws.Range("K139").Value = 8270.143199999999
$ws.Range("L139").Value = 11499
$ws.Range("M139").Value = -3130.143199999999
$ws.Range("N139").Value = -21779

# Row 140
$ws.Range("H140").Value = 1671832.9
$ws.Range("J140").Value = 6999
$ws.Range("L140").Value = 20997
$ws.Range("N140").Value = -31357

# Row 141
$ws.Range("H141").Value = 4611.8
$ws.Range("I141").Value = 4611.8
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 13835.4
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -8655.400000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 521.1818
$ws.Range("I97").Value = 521.6667
$ws.Range("J97").Value = 519
$ws.Range("K97").Value = 521.6667
$ws.Range("L97").Value = 519
$ws.Range("M97").Value = -25.66669999999999
$ws.Range("N97").Value = -1511

# Row 113
$ws.Range("H113").Value = 1705.75
$ws.Range("I113").Value = 1011
$ws.Range("J113").Value = 1937.3334
$ws.Range("K113").Value = 1011
$ws.Range("L113").Value = 1937.3334
$ws.Range("M113").Value = 1159
$ws.Range("N113").Value = -6277.3334

# Row 126
$ws.Range("H126").Value = 1666.3334
$ws.Range("I126").Value = 1666.3334
$ws.Range("K126").Value = 4999.0002
$ws.Range("M126").Value = -2529.0002

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 38125
$ws.Range("I3").Value = 37500
$ws.Range("K3").Value = 37500
$ws.Range("M3").Value = -37388

# Row 15
$ws.Range("H15").Value = 38125
$ws.Range("I15").Value = 37500
$ws.Range("K15").Value = 37500
$ws.Range("M15").Value = -37330

# Row 68
$ws.Range("H68").Value = 5499.8
$ws.Range("I68").Value = 6124.75
$ws.Range("K68").Value = 6124.75
$ws.Range("M68").Value = -5375.75

# Row 71
$ws.Range("H71").Value = 5499.8
$ws.Range("I71").Value = 6124.75
$ws.Range("K71").Value = 30623.75
$ws.Range("M71").Value = -26879.75

$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 38366
$ws.Range("I45").Value = 26997.666
$ws.Range("J45").Value = 45187
$ws.Range("K45").Value = 26997.666
$ws.Range("L45").Value = 45187
$ws.Range("M45").Value = -26506.666
$ws.Range("N45").Value = -46169

# Row 81
$ws.Range("H81").Value = 2537.7646
$ws.Range("J81").Value = 4420
$ws.Range("L81").Value = 8840
$ws.Range("N81").Value = -10962

# Row 84
$ws.Range("H84").Value = 2537.7646
$ws.Range("J84").Value = 4420
$ws.Range("L84").Value = 44200
$ws.Range("N84").Value = -54808

# Row 136
$ws.Range("H136").Value = 3315.6458
$ws.Range("I136").Value = 3347.1428
$ws.Range("J136").Value = 3271.55
$ws.Range("K136").Value = 10041.4284
$ws.Range("L136").Value = 9814.650000000001
$ws.Range("M136").Value = -7491.428400000001
$ws.Range("N136").Value = -14914.65
